# Auto-generated Excel COM-interop script to update market price data
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 1970
$ws.Range("I32").Value = 793.3333
$ws.Range("J32").Value = 5500
$ws.Range("K32").Value = 793.3333
$ws.Range("L32").Value = 5500
$ws.Range("M32").Value = -467.3333
$ws.Range("N32").Value = -6152

# Row 62 (Leve Item ID 27781)
$ws.Range("H62").Value = 7238.5
$ws.Range("I62").Value = 7235
$ws.Range("K62").Value = 7235
$ws.Range("M62").Value = -6611

# Row 65 (Leve Item ID 27781)
$ws.Range("H65").Value = 7238.5
$ws.Range("I65").Value = 7235
$ws.Range("K65").Value = 36175
$ws.Range("M65").Value = -33055

# Row 93 (Leve Item ID 18043)
$ws.Range("H93").Value = 99799.5
$ws.Range("J93").Value = 99799.5
$ws.Range("L93").Value = 99799.5
$ws.Range("N93").Value = -104791.5

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 298746.47
$ws.Range("J138").Value = 1668839.9
$ws.Range("L138").Value = 5006519.699999999
$ws.Range("N138").Value = -5016799.699999999


$ws = $wb.Worksheets.Item("ARM")
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 3288.3333
$ws.Range("I45").Value = 3057.2104
$ws.Range("K45").Value = 3057.2104
$ws.Range("M45").Value = -2680.2104

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 4998.3687
$ws.Range("I61").Value = 2694.2144
$ws.Range("K61").Value = 2694.2144
$ws.Range("M61").Value = -2482.2144

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 2088.5789
$ws.Range("I74").Value = 1909.9412
$ws.Range("J74").Value = 3607
$ws.Range("K74").Value = 1909.9412
$ws.Range("L74").Value = 3607
$ws.Range("M74").Value = -1035.9412
$ws.Range("N74").Value = -5355

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 2088.5789
$ws.Range("I77").Value = 1909.9412
$ws.Range("J77").Value = 3607
$ws.Range("K77").Value = 9549.706
$ws.Range("L77").Value = 18035
$ws.Range("M77").Value = -5181.706
$ws.Range("N77").Value = -26771

# Row 92 (Leve Item ID 18050)
$ws.Range("H92").Value = 50000
$ws.Range("J92").Value = 50000
$ws.Range("L92").Value = 50000
$ws.Range("N92").Value = -54992

# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 2548.111
$ws.Range("I110").Value = 1655.1666
$ws.Range("K110").Value = 1655.1666
$ws.Range("M110").Value = 389.8334

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 4998.3687
$ws.Range("I136").Value = 2694.2144
$ws.Range("K136").Value = 8082.6432
$ws.Range("M136").Value = -5532.6432


$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 2872.4
$ws.Range("I20").Value = 2689
$ws.Range("K20").Value = 2689
$ws.Range("M20").Value = -2442

# Row 21 (Leve Item ID 19542)
$ws.Range("H21").Value = 5000
$ws.Range("J21").Value = 5000
$ws.Range("L21").Value = 5000
$ws.Range("N21").Value = -5472

# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 5810.5
$ws.Range("I86").Value = 8526.875
$ws.Range("J86").Value = 2188.6667
$ws.Range("K86").Value = 8526.875
$ws.Range("L86").Value = 2188.6667
$ws.Range("M86").Value = -7403.875
$ws.Range("N86").Value = -4434.6667

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 5810.5
$ws.Range("I89").Value = 8526.875
$ws.Range("J89").Value = 2188.6667
$ws.Range("K89").Value = 42634.375
$ws.Range("L89").Value = 10943.3335
$ws.Range("M89").Value = -37018.375
$ws.Range("N89").Value = -22175.3335

# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 27778914
$ws.Range("I107").Value = 35715348
$ws.Range("J107").Value = 1395.5
$ws.Range("K107").Value = 35715348
$ws.Range("L107").Value = 1395.5
$ws.Range("M107").Value = -35713428
$ws.Range("N107").Value = -5235.5

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 5835.9023
$ws.Range("I134").Value = 2239.6365
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 6718.9095
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -4183.9095
$ws.Range("N134").Value = -35070


$ws = $wb.Worksheets.Item("CRP")
# Row 7 (Leve Item ID 5361)
$ws.Range("H7").Value = 174.81818
$ws.Range("I7").Value = 120.70588
$ws.Range("J7").Value = 358.8
$ws.Range("K7").Value = 120.70588
$ws.Range("L7").Value = 358.8
$ws.Range("M7").Value = -7.705879999999993
$ws.Range("N7").Value = -584.8

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2334.2083
$ws.Range("I31").Value = 1410.091
$ws.Range("J31").Value = 12499.5
$ws.Range("K31").Value = 1410.091
$ws.Range("L31").Value = 12499.5
$ws.Range("M31").Value = -1115.091
$ws.Range("N31").Value = -13089.5

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2334.2083
$ws.Range("I34").Value = 1410.091
$ws.Range("J34").Value = 12499.5
$ws.Range("K34").Value = 1410.091
$ws.Range("L34").Value = 12499.5
$ws.Range("M34").Value = -1208.091
$ws.Range("N34").Value = -12903.5

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 3275.125
$ws.Range("I58").Value = 3340.2
$ws.Range("K58").Value = 3340.2
$ws.Range("M58").Value = -3137.2

# Row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

# Row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 3559
$ws.Range("I134").Value = 1846.2667
$ws.Range("J134").Value = 5699.9165
$ws.Range("K134").Value = 5538.800099999999
$ws.Range("L134").Value = 17099.7495
$ws.Range("M134").Value = -3003.800099999999
$ws.Range("N134").Value = -22169.7495

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 3275.125
$ws.Range("I136").Value = 3340.2
$ws.Range("K136").Value = 10020.6
$ws.Range("M136").Value = -7470.599999999999


$ws = $wb.Worksheets.Item("CUL")
# Row 26 (Leve Item ID 4746)
$ws.Range("H26").Value = 256.14285
$ws.Range("I26").Value = 249.5
$ws.Range("J26").Value = 258.8
$ws.Range("K26").Value = 748.5
$ws.Range("L26").Value = 776.4000000000001
$ws.Range("M26").Value = -460.5
$ws.Range("N26").Value = -1352.4

# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 1180
$ws.Range("I132").Value = 1098.75
$ws.Range("J132").Value = 1505
$ws.Range("K132").Value = 9888.75
$ws.Range("L132").Value = 13545
$ws.Range("M132").Value = -7358.75
$ws.Range("N132").Value = -18605


$ws = $wb.Worksheets.Item("GSM")
# Row 20 (Leve Item ID 4095)
$ws.Range("H20").Value = 33333
$ws.Range("J20").Value = 33333
$ws.Range("L20").Value = 33333
$ws.Range("N20").Value = -33823

# Row 24 (Leve Item ID 4431)
$ws.Range("H24").Value = 15555
$ws.Range("I24").Value = 15555
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 15555
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -15382
$ws.Range("N24").ClearContents()

# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 3413.9412
$ws.Range("I70").Value = 2840.7273
$ws.Range("K70").Value = 2840.7273
$ws.Range("M70").Value = -2570.7273

# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 3413.9412
$ws.Range("I73").Value = 2840.7273
$ws.Range("K73").Value = 2840.7273
$ws.Range("M73").Value = -1904.7273

# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 23842.6
$ws.Range("I102").Value = 1210.1714
$ws.Range("J102").Value = 103056.1
$ws.Range("K102").Value = 1210.1714
$ws.Range("L102").Value = 103056.1
$ws.Range("M102").Value = 411.8286000000001
$ws.Range("N102").Value = -106300.1


$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 11770354
$ws.Range("I7").Value = 25002940
$ws.Range("J7").Value = 8055.1113
$ws.Range("K7").Value = 25002940
$ws.Range("L7").Value = 8055.1113
$ws.Range("M7").Value = -25002828
$ws.Range("N7").Value = -8279.1113

# Row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2680
$ws.Range("I22").Value = 1119.1428
$ws.Range("J22").Value = 3460.4285
$ws.Range("K22").Value = 1119.1428
$ws.Range("L22").Value = 3460.4285
$ws.Range("M22").Value = -824.1428000000001
$ws.Range("N22").Value = -4050.4285

# Row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2680
$ws.Range("I27").Value = 1119.1428
$ws.Range("J27").Value = 3460.4285
$ws.Range("K27").Value = 1119.1428
$ws.Range("L27").Value = 3460.4285
$ws.Range("M27").Value = -1012.1428
$ws.Range("N27").Value = -3674.4285

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 11770354
$ws.Range("I126").Value = 25002940
$ws.Range("J126").Value = 8055.1113
$ws.Range("K126").Value = 75008820
$ws.Range("L126").Value = 24165.3339
$ws.Range("M126").Value = -75006350
$ws.Range("N126").Value = -29105.3339

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 2645.758
$ws.Range("I132").Value = 2561.34
$ws.Range("K132").Value = 7684.02
$ws.Range("M132").Value = -5154.02

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 4169.7393
$ws.Range("I136").Value = 3952.5789
$ws.Range("K136").Value = 11857.7367
$ws.Range("M136").Value = -9307.736699999999


$ws = $wb.Worksheets.Item("WVR")
# Row 68 (Leve Item ID 10762)
$ws.Range("H68").Value = 83000
$ws.Range("J68").Value = 83000
$ws.Range("L68").Value = 83000
$ws.Range("N68").Value = -84622

# Row 70 (Leve Item ID 11979)
$ws.Range("H70").Value = 30000
$ws.Range("J70").Value = 30000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30630

# Row 71 (Leve Item ID 10762)
$ws.Range("H71").Value = 83000
$ws.Range("J71").Value = 83000
$ws.Range("L71").Value = 249000
$ws.Range("N71").Value = -257112

# Row 73 (Leve Item ID 11979)
$ws.Range("H73").Value = 30000
$ws.Range("J73").Value = 30000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -32184

# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 579.06665
$ws.Range("I107").Value = 299.4
$ws.Range("K107").Value = 898.1999999999999
$ws.Range("M107").Value = 1021.8

# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 5533.8
$ws.Range("I122").Value = 5545.067
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 16635.201
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -14185.201
$ws.Range("N122").Value = -21400

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 3647.3333
$ws.Range("I126").Value = 4
$ws.Range("J126").Value = 4376
$ws.Range("K126").Value = 12
$ws.Range("L126").Value = 13128
$ws.Range("M126").Value = 2458
$ws.Range("N126").Value = -18068

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 2195.9092
$ws.Range("I136").Value = 2215.5
$ws.Range("K136").Value = 6646.5
$ws.Range("M136").Value = -4096.5

